# Remove all traces of "AddressBook" from the Logic Component Class Diagram
# (Dev Guide). The diagram has a single shape whose text reads
# "AddressBook" / "Parser" on two separate paragraphs; the first
# paragraph's text is renamed to "Inventory" while everything else
# (formatting, the second paragraph, every other shape) is left intact.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)

        if (-not $shape.HasTextFrame) { continue }
        if (-not $shape.TextFrame.HasText) { continue }

        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -notlike "*AddressBook*") { continue }

        # Walk the paragraphs of this shape and rewrite only the one that
        # literally says "AddressBook" (leaving sibling paragraphs, such as
        # the "Parser" line directly under it, completely untouched).
        # NB: TextRange.Text on a single paragraph includes the trailing
        # paragraph-mark character (CR), so trim that off before comparing.
        $paraIndex = 1
        while ($paraIndex -le $tr.Paragraphs().Count) {
            $para = $tr.Paragraphs($paraIndex, 1)
            $paraText = $para.Text.TrimEnd([char]13, [char]10)
            if ($paraText -eq "AddressBook") {
                $para.Text = "Inventory"
            }
            $paraIndex = $paraIndex + 1
        }
    }
}
